# edit.ps1 - reproduces the "resultado_horarios" restructuring of Sheet1:
#   - the DIA label column (old F) shifts to E
#   - the even-row time block (old I:J) shifts to G:H
#   - the odd-row time block (old L:M) shifts to J:K
#   - a new column N gets the filter/help labels
#   - column widths / selection are refreshed to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-CellValue($sheet, $src, $dst) {
    # Relocate a cell's value (numbers stay numbers, shared strings stay
    # strings) and remove the old cell entirely so it does not linger as an
    # empty <c> in the saved sheet.
    $cell = $sheet.Range($src)
    $v = $cell.Value2
    $sheet.Range($dst).Value = $v
    $cell.Clear()
}

# ---- Move the DIA label column F -> E ----
Move-CellValue $ws "F2" "E2"
Move-CellValue $ws "F3" "E3"
Move-CellValue $ws "F4" "E4"
Move-CellValue $ws "F5" "E5"
Move-CellValue $ws "F6" "E6"
Move-CellValue $ws "F7" "E7"

# ---- Move the even-row schedule block I:J -> G:H ----
Move-CellValue $ws "I2" "G2"
Move-CellValue $ws "J2" "H2"
Move-CellValue $ws "I4" "G4"
Move-CellValue $ws "J4" "H4"
Move-CellValue $ws "I6" "G6"
Move-CellValue $ws "J6" "H6"
Move-CellValue $ws "I8" "G8"
Move-CellValue $ws "J8" "H8"
Move-CellValue $ws "I10" "G10"
Move-CellValue $ws "J10" "H10"
Move-CellValue $ws "I12" "G12"
Move-CellValue $ws "J12" "H12"
Move-CellValue $ws "I14" "G14"
Move-CellValue $ws "J14" "H14"
Move-CellValue $ws "I16" "G16"
Move-CellValue $ws "J16" "H16"
Move-CellValue $ws "I18" "G18"
Move-CellValue $ws "J18" "H18"
Move-CellValue $ws "I20" "G20"
Move-CellValue $ws "J20" "H20"

# ---- Move the odd-row schedule block L:M -> J:K ----
Move-CellValue $ws "L3" "J3"
Move-CellValue $ws "M3" "K3"
Move-CellValue $ws "L5" "J5"
Move-CellValue $ws "M5" "K5"
Move-CellValue $ws "L7" "J7"
Move-CellValue $ws "M7" "K7"
Move-CellValue $ws "L9" "J9"
Move-CellValue $ws "M9" "K9"
Move-CellValue $ws "L11" "J11"
Move-CellValue $ws "M11" "K11"
Move-CellValue $ws "L13" "J13"
Move-CellValue $ws "M13" "K13"
Move-CellValue $ws "L15" "J15"
Move-CellValue $ws "M15" "K15"
Move-CellValue $ws "L17" "J17"
Move-CellValue $ws "M17" "K17"
Move-CellValue $ws "L19" "J19"
Move-CellValue $ws "M19" "K19"

# ---- New filter/description labels in column N ----
# (written in first-use order so the shared string table matches the
# canonical workbook: Ordenar.. (N6) is interned before Cupo.. (N5))
$ws.Range("N2").Value = "FILTROS DEL CHORARIO"
$ws.Range("N3").Value = "Centro de Idiomas Rocherau"
$ws.Range("N4").Value = "Sede Principal Bogotá"
$ws.Range("N6").Value = "Ordenar por NRC (curso)"
$ws.Range("N5").Value = "Cupo 35 o 45"
$ws.Range("N7").Value = "Para saber el día se concatena todos las columnas  DIA"

# ---- Column widths for the new layout ----
# (ColumnWidth is expressed in characters; the values below are the inputs
# that reproduce the stored widths of the canonical workbook as closely as
# the engine's width quantization allows)
$ws.Columns("B:C").ColumnWidth = 4.25
$ws.Columns("D").ColumnWidth = 2.4166666666666665
$ws.Columns("E").ColumnWidth = 1.9166666666666665
$ws.Columns("F").ColumnWidth = 2.9166666666666665
$ws.Columns("G:H").ColumnWidth = 4.25
$ws.Columns("J:K").ColumnWidth = 4.25
$ws.Columns("N").ColumnWidth = 21.75

# ---- Selection moves to O17 ----
$ws.Range("O17").Select()
